$d = $word.ActiveDocument
try {
  Write-Output ($d.Range().WordOpenXML.Substring(0,200))
} catch {
  Write-Output ("error: " + $_.Exception.Message)
}
